$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column L (correct_ans) abbreviations to full words:
#   b -> center, y -> left, r -> right
$ws.Range("L2").Value = "center"
$ws.Range("L3").Value = "left"
$ws.Range("L4").Value = "right"
$ws.Range("L5").Value = "right"
$ws.Range("L6").Value = "center"
$ws.Range("L7").Value = "left"
$ws.Range("L8").Value = "center"
$ws.Range("L9").Value = "left"
$ws.Range("L10").Value = "right"
$ws.Range("L11").Value = "right"
$ws.Range("L12").Value = "left"
$ws.Range("L13").Value = "center"
$ws.Range("L14").Value = "center"
$ws.Range("L15").Value = "left"
$ws.Range("L16").Value = "right"
$ws.Range("L17").Value = "center"
$ws.Range("L18").Value = "right"
$ws.Range("L19").Value = "left"
$ws.Range("L20").Value = "center"
$ws.Range("L21").Value = "right"
$ws.Range("L22").Value = "left"
$ws.Range("L23").Value = "center"
$ws.Range("L24").Value = "right"
$ws.Range("L25").Value = "left"
$ws.Range("L26").Value = "left"
$ws.Range("L27").Value = "right"
$ws.Range("L28").Value = "center"
$ws.Range("L29").Value = "left"
$ws.Range("L30").Value = "center"
$ws.Range("L31").Value = "right"
$ws.Range("L32").Value = "left"
$ws.Range("L33").Value = "right"
$ws.Range("L34").Value = "center"
$ws.Range("L35").Value = "right"
$ws.Range("L36").Value = "center"
$ws.Range("L37").Value = "left"
$ws.Range("L38").Value = "center"
$ws.Range("L39").Value = "right"
$ws.Range("L40").Value = "left"
$ws.Range("L41").Value = "left"
$ws.Range("L42").Value = "right"
$ws.Range("L43").Value = "center"
$ws.Range("L44").Value = "center"
$ws.Range("L45").Value = "left"
$ws.Range("L46").Value = "right"
$ws.Range("L47").Value = "left"
$ws.Range("L48").Value = "center"
$ws.Range("L49").Value = "right"
$ws.Range("L50").Value = "center"
$ws.Range("L51").Value = "left"
$ws.Range("L52").Value = "right"
$ws.Range("L53").Value = "left"
$ws.Range("L54").Value = "right"
$ws.Range("L55").Value = "center"
$ws.Range("L56").Value = "center"
$ws.Range("L57").Value = "left"
$ws.Range("L58").Value = "right"
$ws.Range("L59").Value = "center"
$ws.Range("L60").Value = "left"
$ws.Range("L61").Value = "right"
$ws.Range("L62").Value = "left"
$ws.Range("L63").Value = "right"
$ws.Range("L64").Value = "center"
$ws.Range("L65").Value = "center"
$ws.Range("L66").Value = "right"
$ws.Range("L67").Value = "left"
$ws.Range("L68").Value = "left"
$ws.Range("L69").Value = "right"
$ws.Range("L70").Value = "center"
$ws.Range("L71").Value = "right"
$ws.Range("L72").Value = "left"
$ws.Range("L73").Value = "center"
$ws.Range("L74").Value = "right"
$ws.Range("L75").Value = "center"
$ws.Range("L76").Value = "left"
$ws.Range("L77").Value = "left"
$ws.Range("L78").Value = "right"
$ws.Range("L79").Value = "center"
$ws.Range("L80").Value = "center"
$ws.Range("L81").Value = "right"
$ws.Range("L82").Value = "left"
$ws.Range("L83").Value = "right"
$ws.Range("L84").Value = "center"
$ws.Range("L85").Value = "left"
$ws.Range("L86").Value = "right"
$ws.Range("L87").Value = "center"
$ws.Range("L88").Value = "left"
$ws.Range("L89").Value = "right"
$ws.Range("L90").Value = "center"
$ws.Range("L91").Value = "left"
$ws.Range("L92").Value = "right"
$ws.Range("L93").Value = "left"
$ws.Range("L94").Value = "center"
$ws.Range("L95").Value = "right"
$ws.Range("L96").Value = "left"
$ws.Range("L97").Value = "center"
$ws.Range("L98").Value = "center"
$ws.Range("L99").Value = "left"
$ws.Range("L100").Value = "right"
$ws.Range("L101").Value = "center"
$ws.Range("L102").Value = "right"
$ws.Range("L103").Value = "left"
$ws.Range("L104").Value = "left"
$ws.Range("L105").Value = "right"
$ws.Range("L106").Value = "center"
$ws.Range("L107").Value = "left"
$ws.Range("L108").Value = "center"
$ws.Range("L109").Value = "right"
$ws.Range("L110").Value = "center"
$ws.Range("L111").Value = "right"
$ws.Range("L112").Value = "left"
$ws.Range("L113").Value = "center"
$ws.Range("L114").Value = "left"
$ws.Range("L115").Value = "right"
$ws.Range("L116").Value = "left"
$ws.Range("L117").Value = "center"
$ws.Range("L118").Value = "right"
$ws.Range("L119").Value = "center"
$ws.Range("L120").Value = "left"
$ws.Range("L121").Value = "right"
$ws.Range("L122").Value = "center"
$ws.Range("L123").Value = "left"
$ws.Range("L124").Value = "right"
$ws.Range("L125").Value = "center"
$ws.Range("L126").Value = "left"
$ws.Range("L127").Value = "right"
$ws.Range("L128").Value = "center"
$ws.Range("L129").Value = "right"
$ws.Range("L130").Value = "left"
$ws.Range("L131").Value = "center"
$ws.Range("L132").Value = "right"
$ws.Range("L133").Value = "left"
$ws.Range("L134").Value = "center"
$ws.Range("L135").Value = "right"
$ws.Range("L136").Value = "left"
$ws.Range("L137").Value = "center"
$ws.Range("L138").Value = "left"
$ws.Range("L139").Value = "right"
$ws.Range("L140").Value = "right"
$ws.Range("L141").Value = "left"
$ws.Range("L142").Value = "center"
$ws.Range("L143").Value = "center"
$ws.Range("L144").Value = "left"
$ws.Range("L145").Value = "right"
$ws.Range("L146").Value = "right"
$ws.Range("L147").Value = "center"
$ws.Range("L148").Value = "left"
$ws.Range("L149").Value = "right"
$ws.Range("L150").Value = "left"
$ws.Range("L151").Value = "center"
$ws.Range("L152").Value = "right"
$ws.Range("L153").Value = "left"
$ws.Range("L154").Value = "center"
$ws.Range("L155").Value = "right"
$ws.Range("L156").Value = "center"
$ws.Range("L157").Value = "left"
$ws.Range("L158").Value = "center"
$ws.Range("L159").Value = "left"
$ws.Range("L160").Value = "right"
$ws.Range("L161").Value = "center"
$ws.Range("L162").Value = "left"
$ws.Range("L163").Value = "right"
$ws.Range("L164").Value = "left"
$ws.Range("L165").Value = "center"
$ws.Range("L166").Value = "right"
$ws.Range("L167").Value = "left"
$ws.Range("L168").Value = "right"
$ws.Range("L169").Value = "center"
$ws.Range("L170").Value = "right"
$ws.Range("L171").Value = "left"
$ws.Range("L172").Value = "center"
$ws.Range("L173").Value = "left"
$ws.Range("L174").Value = "center"
$ws.Range("L175").Value = "right"
$ws.Range("L176").Value = "right"
$ws.Range("L177").Value = "left"
$ws.Range("L178").Value = "center"
$ws.Range("L179").Value = "right"
$ws.Range("L180").Value = "center"
$ws.Range("L181").Value = "left"
$ws.Range("L182").Value = "left"
$ws.Range("L183").Value = "center"
$ws.Range("L184").Value = "right"
$ws.Range("L185").Value = "left"
$ws.Range("L186").Value = "center"
$ws.Range("L187").Value = "right"
$ws.Range("L188").Value = "right"
$ws.Range("L189").Value = "left"
$ws.Range("L190").Value = "center"
$ws.Range("L191").Value = "center"
$ws.Range("L192").Value = "right"
$ws.Range("L193").Value = "left"
$ws.Range("L194").Value = "center"
$ws.Range("L195").Value = "right"
$ws.Range("L196").Value = "left"
$ws.Range("L197").Value = "left"
$ws.Range("L198").Value = "right"
$ws.Range("L199").Value = "center"
$ws.Range("L200").Value = "left"
$ws.Range("L201").Value = "right"
$ws.Range("L202").Value = "center"
$ws.Range("L203").Value = "left"
$ws.Range("L204").Value = "right"
$ws.Range("L205").Value = "center"
$ws.Range("L206").Value = "left"
$ws.Range("L207").Value = "center"
$ws.Range("L208").Value = "right"
$ws.Range("L209").Value = "left"
$ws.Range("L210").Value = "right"
$ws.Range("L211").Value = "center"
$ws.Range("L212").Value = "left"
$ws.Range("L213").Value = "center"
$ws.Range("L214").Value = "right"
$ws.Range("L215").Value = "right"
$ws.Range("L216").Value = "left"
$ws.Range("L217").Value = "center"
$ws.Range("L218").Value = "center"
$ws.Range("L219").Value = "left"
$ws.Range("L220").Value = "right"
$ws.Range("L221").Value = "right"
$ws.Range("L222").Value = "left"
$ws.Range("L223").Value = "center"
$ws.Range("L224").Value = "left"
$ws.Range("L225").Value = "right"
$ws.Range("L226").Value = "center"
$ws.Range("L227").Value = "left"
$ws.Range("L228").Value = "right"
$ws.Range("L229").Value = "center"
$ws.Range("L230").Value = "right"
$ws.Range("L231").Value = "left"
$ws.Range("L232").Value = "center"
$ws.Range("L233").Value = "left"
$ws.Range("L234").Value = "right"
$ws.Range("L235").Value = "center"
$ws.Range("L236").Value = "right"
$ws.Range("L237").Value = "center"
$ws.Range("L238").Value = "left"
$ws.Range("L239").Value = "left"
$ws.Range("L240").Value = "center"
$ws.Range("L241").Value = "right"
$ws.Range("L242").Value = "center"
$ws.Range("L243").Value = "right"
$ws.Range("L244").Value = "left"
$ws.Range("L245").Value = "right"
$ws.Range("L246").Value = "left"
$ws.Range("L247").Value = "center"
$ws.Range("L248").Value = "right"
$ws.Range("L249").Value = "left"
$ws.Range("L250").Value = "center"
$ws.Range("L251").Value = "right"
$ws.Range("L252").Value = "center"
$ws.Range("L253").Value = "left"
$ws.Range("L254").Value = "right"
$ws.Range("L255").Value = "left"
$ws.Range("L256").Value = "center"
$ws.Range("L257").Value = "left"
$ws.Range("L258").Value = "right"
$ws.Range("L259").Value = "center"
$ws.Range("L260").Value = "center"
$ws.Range("L261").Value = "left"
$ws.Range("L262").Value = "right"
$ws.Range("L263").Value = "left"
$ws.Range("L264").Value = "right"
$ws.Range("L265").Value = "center"
$ws.Range("L266").Value = "center"
$ws.Range("L267").Value = "right"
$ws.Range("L268").Value = "left"
$ws.Range("L269").Value = "center"
$ws.Range("L270").Value = "right"
$ws.Range("L271").Value = "left"
$ws.Range("L272").Value = "right"
$ws.Range("L273").Value = "left"
$ws.Range("L274").Value = "center"
$ws.Range("L275").Value = "center"
$ws.Range("L276").Value = "right"
$ws.Range("L277").Value = "left"
$ws.Range("L278").Value = "right"
$ws.Range("L279").Value = "center"
$ws.Range("L280").Value = "left"
$ws.Range("L281").Value = "left"
$ws.Range("L282").Value = "right"
$ws.Range("L283").Value = "center"
$ws.Range("L284").Value = "center"
$ws.Range("L285").Value = "left"
$ws.Range("L286").Value = "right"
$ws.Range("L287").Value = "center"
$ws.Range("L288").Value = "left"
$ws.Range("L289").Value = "right"
$ws.Range("L290").Value = "left"
$ws.Range("L291").Value = "right"
$ws.Range("L292").Value = "center"
$ws.Range("L293").Value = "left"
$ws.Range("L294").Value = "center"
$ws.Range("L295").Value = "right"
$ws.Range("L296").Value = "center"
$ws.Range("L297").Value = "left"
$ws.Range("L298").Value = "right"
$ws.Range("L299").Value = "center"
$ws.Range("L300").Value = "right"
$ws.Range("L301").Value = "left"
$ws.Range("L302").Value = "right"
$ws.Range("L303").Value = "center"
$ws.Range("L304").Value = "left"
$ws.Range("L305").Value = "left"
$ws.Range("L306").Value = "right"
$ws.Range("L307").Value = "center"
$ws.Range("L308").Value = "left"
$ws.Range("L309").Value = "center"
$ws.Range("L310").Value = "right"
$ws.Range("L311").Value = "right"
$ws.Range("L312").Value = "left"
$ws.Range("L313").Value = "center"
$ws.Range("L314").Value = "left"
$ws.Range("L315").Value = "right"
$ws.Range("L316").Value = "center"
$ws.Range("L317").Value = "right"
$ws.Range("L318").Value = "left"
$ws.Range("L319").Value = "center"
$ws.Range("L320").Value = "left"
$ws.Range("L321").Value = "right"
$ws.Range("L322").Value = "center"
$ws.Range("L323").Value = "center"
$ws.Range("L324").Value = "left"
$ws.Range("L325").Value = "right"
$ws.Range("L326").Value = "right"
$ws.Range("L327").Value = "left"
$ws.Range("L328").Value = "center"
$ws.Range("L329").Value = "left"
$ws.Range("L330").Value = "center"
$ws.Range("L331").Value = "right"
$ws.Range("L332").Value = "center"
$ws.Range("L333").Value = "left"
$ws.Range("L334").Value = "right"
$ws.Range("L335").Value = "center"
$ws.Range("L336").Value = "right"
$ws.Range("L337").Value = "left"
$ws.Range("L338").Value = "right"
$ws.Range("L339").Value = "left"
$ws.Range("L340").Value = "center"
$ws.Range("L341").Value = "left"
$ws.Range("L342").Value = "center"
$ws.Range("L343").Value = "right"
$ws.Range("L344").Value = "center"
$ws.Range("L345").Value = "left"
$ws.Range("L346").Value = "right"
$ws.Range("L347").Value = "right"
$ws.Range("L348").Value = "left"
$ws.Range("L349").Value = "center"
$ws.Range("L350").Value = "center"
$ws.Range("L351").Value = "left"
$ws.Range("L352").Value = "right"
$ws.Range("L353").Value = "left"
$ws.Range("L354").Value = "right"
$ws.Range("L355").Value = "center"
$ws.Range("L356").Value = "left"
$ws.Range("L357").Value = "center"
$ws.Range("L358").Value = "right"
$ws.Range("L359").Value = "center"
$ws.Range("L360").Value = "right"
$ws.Range("L361").Value = "left"

# Update image folder references from "face" to "book"
$ws.Range("D21").Value = "book//book_06.jpg"
$ws.Range("D34").Value = "book//book_26.jpg"
$ws.Range("D65").Value = "book//book_20.jpg"
$ws.Range("D81").Value = "book//book_05.jpg"
$ws.Range("D96").Value = "book//book_30.jpg"
$ws.Range("D103").Value = "book//book_16.jpg"
$ws.Range("D111").Value = "book//book_30.jpg"
$ws.Range("D157").Value = "book//book_09.jpg"
$ws.Range("D160").Value = "book//book_02.jpg"
$ws.Range("D165").Value = "book//book_16.jpg"
$ws.Range("D190").Value = "book//book_08.jpg"
$ws.Range("D195").Value = "book//book_25.jpg"
$ws.Range("D210").Value = "book//book_39.jpg"
$ws.Range("D289").Value = "book//book_21.jpg"
$ws.Range("D290").Value = "book//book_20.jpg"
$ws.Range("B302").Value = "book//book_08.jpg"
$ws.Range("A303").Value = "book//book_08.jpg"
$ws.Range("B308").Value = "book//book_25.jpg"
$ws.Range("A309").Value = "book//book_25.jpg"
$ws.Range("C311").Value = "book//book_25.jpg"
$ws.Range("B314").Value = "book//book_19.jpg"
$ws.Range("A315").Value = "book//book_19.jpg"
$ws.Range("C319").Value = "book//book_19.jpg"
$ws.Range("B320").Value = "book//book_15.jpg"
$ws.Range("A321").Value = "book//book_15.jpg"
$ws.Range("C324").Value = "book//book_15.jpg"
$ws.Range("B326").Value = "book//book_18.jpg"
$ws.Range("A327").Value = "book//book_18.jpg"
$ws.Range("C329").Value = "book//book_18.jpg"
$ws.Range("B332").Value = "book//book_27.jpg"
$ws.Range("A333").Value = "book//book_27.jpg"
$ws.Range("B338").Value = "book//book_21.jpg"
$ws.Range("A339").Value = "book//book_21.jpg"
$ws.Range("B344").Value = "book//book_12.jpg"
$ws.Range("A345").Value = "book//book_12.jpg"
$ws.Range("C348").Value = "book//book_12.jpg"
$ws.Range("C349").Value = "book//book_12.jpg"
$ws.Range("B350").Value = "book//book_22.jpg"
$ws.Range("A351").Value = "book//book_22.jpg"
$ws.Range("B356").Value = "book//book_32.jpg"
$ws.Range("A357").Value = "book//book_32.jpg"
$ws.Range("C360").Value = "book//book_32.jpg"
